# Apply the two content changes described in the commit message:
#   1. "Volume P-20" -> "Volume D-20"
#   2. "Lecture Notes in Informatics (LNI) - Proceedings"
#        -> "Lecture Notes in Informatics (LNI) - Dissertations"
#
# Both target strings each occur exactly once in the document and sit
# entirely inside a single existing run, so a plain Find/Replace keeps the
# run's character formatting (w:rPr, incl. w:lang) intact.

$d = $word.ActiveDocument

# 1) Sub-title on the "Volume ..." header: Proceedings volume "P" becomes
#    the Dissertations volume "D". (This also naturally drops the stray
#    "_GoBack" bookmark that used to sit around the lone "P" run, because
#    the whole matched range -- including the bookmark -- gets replaced.)
$d.Content.Find.Execute("Volume P-20", $true, $false, $false, $false, $false, $true, 1, $false, "Volume D-20", 2) | Out-Null

# 2) Strang sub-title: "... Proceedings" -> "... Dissertations"
$d.Content.Find.Execute("Lecture Notes in Informatics (LNI) - Proceedings", $true, $false, $false, $false, $false, $true, 1, $false, "Lecture Notes in Informatics (LNI) - Dissertations", 2) | Out-Null
